$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.095.00"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -3.39%  "
$ws.Range("D3").Value = "'3.331.07"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -5.55%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'550.56"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.07%  "
$ws.Range("D6").Value = "'172.54"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.04%  "
$ws.Range("E7").Value = "  -3.43%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'3.324.35"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.50%  "
$ws.Range("E10").Value = "  -2.85%  "
$ws.Range("D11").Value = "'0.161"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.33%  "
$ws.Range("D12").Value = "'53.16"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.63%  "
$ws.Range("E13").Value = "  -1.98%  "
$ws.Range("D14").Value = "'9.00"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.37%  "
$ws.Range("D15").Value = "'3.859.50"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.53%  "
$ws.Range("D16").Value = "'18.25"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'3.337.67"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.96%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.117"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.85%  "
$ws.Range("D19").Value = "'11.73"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.78%  "
$ws.Range("D20").Value = "'63.934.18"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.58%  "
$ws.Range("D21").Value = "'0.972"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.84%  "
$ws.Range("D22").Value = "'424.08"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.98%  "
$ws.Range("D23").Value = "'4.74"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +10.97%  "
$ws.Range("E24").Value = "  -2.71%  "
$ws.Range("D25").Value = "'83.97"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.14%  "
$ws.Range("D26").Value = "'13.21"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.02%  "
$ws.Range("D27").Value = "'10.61"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.38%  "
$ws.Range("D28").Value = "'2.81"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.22%  "
$ws.Range("D29").Value = "'8.57"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.76%  "
$ws.Range("D30").Value = "'29.63"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.12%  "
$ws.Range("E31").Value = "  +1.65%  "
$ws.Range("D32").Value = "'593.73"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -7.55%  "
$ws.Range("D33").Value = "'11.39"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.11%  "
$ws.Range("E34").Value = "  -4.35%  "
$ws.Range("D35").Value = "'58.17"
$ws.Range("D35").ClearFormats()
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "'0.142"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -9.47%  "
$ws.Range("D38").Value = "'3.49"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("D39").Value = "'35.35"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.74%  "
$ws.Range("E40").Value = "  -8.05%  "
$ws.Range("E41").Value = "  -5.02%  "
$ws.Range("D42").Value = "'3.099.99"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.09%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").Value = "'2.79"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.06%  "
$ws.Range("D45").Value = "'3.19"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.88%  "
$ws.Range("E46").Value = "  -3.79%  "
$ws.Range("E47").Value = "  -4.36%  "
$ws.Range("E48").Value = "  -3.34%  "
$ws.Range("E49").Value = "  -4.72%  "
$ws.Range("D50").Value = "'8.18"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -6.28%  "
$ws.Range("D51").Value = "'132.95"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.79%  "
